$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Jelinek-Mercer smoothing row (row 7) with new values
$ws.Range("B7").Value = 0.236
$ws.Range("C7").Value = 0.368
$ws.Range("D7").Value = 0.3093

# Move the active cell selection from B6 to B7
$ws.Range("B7").Select()
